$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Mayo de 2020 a las 22:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1523544
$ws.Range("C4").Value = 15771
$ws.Range("D4").Value = 342690
$ws.Range("E4").Value = 1089962
$ws.Range("G4").Value = 779
$ws.Range("H4").Value = 90892

# Row 8 - Brasil
$ws.Range("B8").Value = 235331
$ws.Range("C8").Value = 2189
$ws.Range("E8").Value = 129935
$ws.Range("G8").Value = 91
$ws.Range("H8").Value = 15724

# Row 15 - Peru
$ws.Range("B15").Value = 92273
$ws.Range("C15").Value = 3732
$ws.Range("E15").Value = 61353
$ws.Range("G15").Value = 125
$ws.Range("H15").Value = 2648

# Row 24 - Ecuador
$ws.Range("B24").Value = 33182
$ws.Range("C24").Value = 419
$ws.Range("E24").Value = 27013
$ws.Range("G24").Value = 48
$ws.Range("H24").Value = 2736

# Row 75 - Uzbekistan
$ws.Range("B75").Value = 2753
$ws.Range("C75").Value = 15
$ws.Range("E75").Value = 494
